$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-converted to a number
# by Excel (single decimal point, no thousands separators) are forced back
# to Text format first so they stay strings, matching the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "41.523.85"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "2.464.55"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("E4").Value = "  -0.63%  "

$ws.Range("D5").Value = "313.83"
$ws.Range("E5").Value = "  +0.11%  "

$ws.Range("D6").Value = "90.84"
$ws.Range("E6").Value = "  -1.64%  "

$ws.Range("D7").Value = "0.549"
$ws.Range("E7").Value = "  +0.91%  "

$ws.Range("E8").Value = "  -0.61%  "

$ws.Range("E9").Value = "  +3.81%  "

$ws.Range("D10").Value = "32.44"
$ws.Range("E10").Value = "  -1.37%  "

$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  +1.61%  "

$ws.Range("E12").Value = "  +0.62%  "

$ws.Range("D13").Value = "2.843.08"
$ws.Range("E13").Value = "  -0.43%  "

$ws.Range("D14").Value = "6.84"
$ws.Range("E14").Value = "  -0.10%  "

$ws.Range("E15").Value = "  +3.37%  "

$ws.Range("D16").Value = "2.515.60"
$ws.Range("E16").Value = "  +2.87%  "

$ws.Range("D17").Value = "0.775"
$ws.Range("E17").Value = "  -1.09%  "

$ws.Range("D18").Value = "41.513.91"
$ws.Range("E18").Value = "  +0.32%  "

$ws.Range("D19").Value = "6.49"
$ws.Range("E19").Value = "  +3.35%  "

$ws.Range("D20").Value = "0.0₃0939"
$ws.Range("E20").Value = "  +2.04%  "

$ws.Range("D21").Value = "71.08"
$ws.Range("E21").Value = "  +1.16%  "

$ws.Range("D22").Value = "11.11"
$ws.Range("E22").Value = "  +0.92%  "

$ws.Range("D23").Value = "237.98"
$ws.Range("E23").Value = "  +1.26%  "

$ws.Range("E24").Value = "  -0.70%  "

$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("E26").Value = "  +1.49%  "

$ws.Range("D27").Value = "24.51"
$ws.Range("E27").Value = "  +2.15%  "

$ws.Range("E28").Value = "  +0.18%  "

$ws.Range("D29").Value = "9.64"
$ws.Range("E29").Value = "  -0.82%  "

$ws.Range("D30").Value = "35.21"
$ws.Range("E30").Value = "  -2.63%  "

$ws.Range("D31").Value = "156.63"
$ws.Range("E31").Value = "  +3.02%  "

$ws.Range("D32").Value = "5.42"
$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("E33").Value = "  +0.94%  "

$ws.Range("D34").Value = "0.0756"
$ws.Range("E34").Value = "  +0.87%  "

$ws.Range("D35").Value = "17.22"
$ws.Range("E35").Value = "  -2.03%  "

$ws.Range("D36").Value = "2.35"
$ws.Range("E36").Value = "  -8.22%  "

$ws.Range("E37").Value = "  -4.26%  "

$ws.Range("E38").Value = "  +1.32%  "

$ws.Range("D39").Value = "0.102"
$ws.Range("E39").Value = "  +3.04%  "

$ws.Range("E40").Value = "  -3.64%  "

$ws.Range("E41").Value = "  -1.40%  "

$ws.Range("E42").Value = "  -0.80%  "

$ws.Range("D43").Value = "1.942.47"
$ws.Range("E43").Value = "  -0.92%  "

$ws.Range("E44").Value = "  +0.19%  "

$ws.Range("D45").Value = "18.62"
$ws.Range("E45").Value = "  -4.02%  "

$ws.Range("D46").Value = "2.89"
$ws.Range("E46").Value = "  -1.94%  "

$ws.Range("D47").Value = "9.02"
$ws.Range("E47").Value = "  +3.74%  "

$ws.Range("D48").Value = "2.703.02"
$ws.Range("E48").Value = "  -0.31%  "

$ws.Range("D49").Value = "96.79"
$ws.Range("E49").Value = "  +1.23%  "

$ws.Range("D50").Value = "67.04"
$ws.Range("E50").Value = "  -2.04%  "

$ws.Range("E51").Value = "  -2.40%  "
